$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.768.63'
$ws.Range('E2').Value = '  +2.15%  '

$ws.Range('D3').Value = '2.117.12'
$ws.Range('E3').Value = '  +10.38%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '334.81'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.19%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5242'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.59%  '

$ws.Range('E8').Value = '  +8.28%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09088'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.95%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.35'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.81%  '

$ws.Range('E11').Value = '  +6.86%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.30'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.41%  '

$ws.Range('D13').Value = '2.115.67'
$ws.Range('E13').Value = '  +10.49%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.779'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.54%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.859'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '98.13'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.08%  '

$ws.Range('E17').Value = '  -0.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001140'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06650'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.07%  '

$ws.Range('E20').Value = '  +3.73%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9994'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.20%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.408'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.78%  '

$ws.Range('D23').Value = '30.894.00'
$ws.Range('E23').Value = '  +2.53%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.09'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.30%  '

$ws.Range('D25').Value = '2.363.25'
$ws.Range('E25').Value = '  +10.78%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.253'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.73%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.01'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.14%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.561'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.03%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.46'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.54%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.76'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.69%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.184'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.39%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1074'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.78%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.269'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.29%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.939'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.540'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +28.56%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02602'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.91%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.603'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.13%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.646'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.11%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06781'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.32%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.80'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.02%  '

$ws.Range('E41').Value = '  +5.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6842'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.260'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.26'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.48%  '

$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6446'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.33%  '

$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9990'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.267'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.51%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.677'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.39%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.289'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.35%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.26'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07089'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.86%  '
